{"js": "// Update the division-problem worksheet table: replace the text of each\n// non-blank problem cell (in reading order: row by row, left to right)\n// with the corresponding new problem. We write into the paragraph's\n// Range (not the cell body as a whole) via InsertLocation.Replace so the\n// existing run/paragraph formatting (font, size, alignment) is preserved.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\n// Old problems, in reading order, used only to sanity-check we are\n// targeting the right cells before overwriting them.\nconst oldValues = [\n  \"709\u00f74=\", \"983\u00f79=\", \"221\u00f77=\", \"287\u00f75=\", \"532\u00f76=\",\n  \"863\u00f74=\", \"792\u00f73=\", \"403\u00f79=\", \"593\u00f79=\", \"912\u00f79=\",\n  \"317\u00f72=\", \"436\u00f76=\", \"147\u00f75=\", \"933\u00f77=\", \"226\u00f74=\",\n  \"335\u00f79=\", \"449\u00f78=\", \"789\u00f79=\", \"597\u00f74=\", \"799\u00f72=\",\n  \"509\u00f72=\", \"307\u00f74=\", \"572\u00f74=\", \"729\u00f77=\", \"666\u00f78=\"\n];\n\n// New problems, same order.\nconst newValues = [\n  \"611\u00f73=\", \"527\u00f77=\", \"652\u00f72=\", \"468\u00f73=\", \"142\u00f75=\",\n  \"502\u00f74=\", \"561\u00f79=\", \"375\u00f79=\", \"773\u00f74=\", \"981\u00f75=\",\n  \"219\u00f76=\", \"301\u00f78=\", \"377\u00f76=\", \"515\u00f79=\", \"556\u00f78=\",\n  \"233\u00f78=\", \"280\u00f73=\", \"362\u00f72=\", \"123\u00f79=\", \"125\u00f75=\",\n  \"820\u00f73=\", \"221\u00f78=\", \"440\u00f76=\", \"509\u00f72=\", \"738\u00f77=\"\n];\n\nconst rowCount = table.rowCount;\nconst colCount = table.values[0].length;\n\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cellText = table.values[r][c];\n    if (cellText !== \"\") {\n      if (idx < oldValues.length && cellText === oldValues[idx]) {\n        const cell = table.getCell(r, c);\n        const range = cell.body.getRange();\n        range.insertText(newValues[idx], Word.InsertLocation.replace);\n      }\n      idx++;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem worksheet table: replace the text of each\n# non-blank problem cell (in reading order: row by row, left to right)\n# with the corresponding new problem, preserving all existing run/paragraph\n# formatting (font, size, alignment) since we only touch Range.Text.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Old problems, in reading order, used only to sanity-check we are\n# targeting the right cells before overwriting them.\n$oldValues = @(\n  \"709\u00f74=\", \"983\u00f79=\", \"221\u00f77=\", \"287\u00f75=\", \"532\u00f76=\",\n  \"863\u00f74=\", \"792\u00f73=\", \"403\u00f79=\", \"593\u00f79=\", \"912\u00f79=\",\n  \"317\u00f72=\", \"436\u00f76=\", \"147\u00f75=\", \"933\u00f77=\", \"226\u00f74=\",\n  \"335\u00f79=\", \"449\u00f78=\", \"789\u00f79=\", \"597\u00f74=\", \"799\u00f72=\",\n  \"509\u00f72=\", \"307\u00f74=\", \"572\u00f74=\", \"729\u00f77=\", \"666\u00f78=\"\n)\n\n# New problems, same order.\n$newValues = @(\n  \"611\u00f73=\", \"527\u00f77=\", \"652\u00f72=\", \"468\u00f73=\", \"142\u00f75=\",\n  \"502\u00f74=\", \"561\u00f79=\", \"375\u00f79=\", \"773\u00f74=\", \"981\u00f75=\",\n  \"219\u00f76=\", \"301\u00f78=\", \"377\u00f76=\", \"515\u00f79=\", \"556\u00f78=\",\n  \"233\u00f78=\", \"280\u00f73=\", \"362\u00f72=\", \"123\u00f79=\", \"125\u00f75=\",\n  \"820\u00f73=\", \"221\u00f78=\", \"440\u00f76=\", \"509\u00f72=\", \"738\u00f77=\"\n)\n\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cellRange = $cell.Range\n    $cellText = $cellRange.Text\n    # Strip the trailing end-of-cell marker(s) (cell mark \\r\\a) before\n    # comparing against the expected problem text.\n    $trimmed = $cellText.TrimEnd([char]13, [char]7)\n    if ($trimmed -ne \"\") {\n      if ($idx -lt $oldValues.Length -and $trimmed -eq $oldValues[$idx]) {\n        $cellRange.Text = $newValues[$idx]\n      }\n      $idx = $idx + 1\n    }\n  }\n}\n"}
